$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9042177796363831
$ws.Range("B1").Value = 1.906989932060242
$ws.Range("C1").Value = 8.618441581726074
$ws.Range("D1").Value = 1.98227322101593
$ws.Range("E1").Value = 1.451583981513977
